# Apply edits described by the diff to LOQ4023.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextKeepingFormat {
    # Positional parameters: worksheet, target address, format-source address, text
    param($Worksheet, [string]$TargetAddress, [string]$FormatSourceAddress, [string]$Text)
    $cell = $Worksheet.Range($TargetAddress)
    # Force text interpretation so Excel does not auto-convert date-like
    # strings (e.g. "01/01/2022") into date serial numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    # Restore the exact visual formatting (font/alignment/number format)
    # from a known-good neighbor cell so the style index matches the
    # original layout instead of the ad-hoc "@" format just applied.
    $Worksheet.Range($FormatSourceAddress).Copy()
    $cell.PasteSpecial(-4122)
}

# Ativação: date changed from 01/01/2012 to 01/01/2022
Set-TextKeepingFormat $ws "B8" "B9" "01/01/2022"
Set-TextKeepingFormat $ws "C8" "C9" "01/01/2022"

# Objetivos: (Portuguese) text updated
$ws.Range("B10").Value = "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica."
$ws.Range("C10").Value = "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica."

# Objectives: (English) text newly added in row 11
Set-TextKeepingFormat $ws "B11" "B10" "Check the students a practical notion of applications on an industrial scale, processes and products studied in organic chemistry disciplines and also an overview of the most important raw materials in the chemical industry of organic base"
Set-TextKeepingFormat $ws "C11" "C10" "Check the students a practical notion of applications on an industrial scale, processes and products studied in organic chemistry disciplines and also an overview of the most important raw materials in the chemical industry of organic base"

# Programa resumido: (Portuguese) text updated
$ws.Range("B14").Value = "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação."
$ws.Range("C14").Value = "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação."

# Short syllabus: (English) text newly added in row 15
Set-TextKeepingFormat $ws "B15" "B14" "Oil, Natural Gas and Petrochemicals; Fine Chemistry; Organic Unit Processes: Nitration; Esterification; Alkylation and Acylation; Hydrogenation; Sulphonation/Sulfation; Oxidation."
Set-TextKeepingFormat $ws "C15" "C14" "Oil, Natural Gas and Petrochemicals; Fine Chemistry; Organic Unit Processes: Nitration; Esterification; Alkylation and Acylation; Hydrogenation; Sulphonation/Sulfation; Oxidation."

# Programa: (Portuguese) text updated
$ws.Range("B16").Value = "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação."
$ws.Range("C16").Value = "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação."

# Syllabus: (English) text newly added in row 17
Set-TextKeepingFormat $ws "B17" "B16" "1- Oil, Natural Gas and Petrochemicals 2- Fine Chemistry: Characteristics, Fine Chemistry vs. Basic Chemistry, Main Segments (Pesticides, Drugs, Catalysts, Dyes and Pigments, Specialties); 3- Organic Unit Processes: 3.1- Nitration; 3.2- Esterification; 3.3- Alkylation and Acylation; 3.4- Hydrogenation; 3.5- Sulphonation/Sulfation; 3.6- Oxidation."
Set-TextKeepingFormat $ws "C17" "C16" "1- Oil, Natural Gas and Petrochemicals 2- Fine Chemistry: Characteristics, Fine Chemistry vs. Basic Chemistry, Main Segments (Pesticides, Drugs, Catalysts, Dyes and Pigments, Specialties); 3- Organic Unit Processes: 3.1- Nitration; 3.2- Esterification; 3.3- Alkylation and Acylation; 3.4- Hydrogenation; 3.5- Sulphonation/Sulfation; 3.6- Oxidation."

# Método: text updated
$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos."

# Critério: text updated
$ws.Range("B20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula."

# Norma de recuperação: text updated
$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."

# Bibliografia: text updated
$ws.Range("B22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."
$ws.Range("C22").Value = "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/."

Write-Host "Edit applied successfully"